# Add the 2018-08-14 and 2018-11-06 elections as two new rows (101, 102)
# at the bottom of the office_table sheet, following the existing pattern
# used by the other rows in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 101 : 2018-08-14 (Partisan Primary) ---------------------------
$ws.Cells.Item(101, 1).Value = 1878          # id
$ws.Cells.Item(101, 2).NumberFormat = "@"    # date column is stored as text
$ws.Cells.Item(101, 2).Value = "2018-08-14"  # date
$ws.Cells.Item(101, 4).Value = "P"           # primary
$ws.Cells.Item(101, 8).Value = "X"           # Senate
$ws.Cells.Item(101, 9).Value = "X"           # House
$ws.Cells.Item(101, 10).Value = "X"          # Governor
$ws.Cells.Item(101, 11).Value = "X"          # Lt Gov
$ws.Cells.Item(101, 12).Value = "X"          # Atty General
$ws.Cells.Item(101, 13).Value = "X"          # Sec of St
$ws.Cells.Item(101, 14).Value = "X"          # St Treasurer
$ws.Cells.Item(101, 16).Value = "X"          # St Senate
$ws.Cells.Item(101, 17).Value = "X"          # St Assembly

# --- Row 102 : 2018-11-06 (General election) ----------------------------
$ws.Cells.Item(102, 1).Value = 1886          # id
$ws.Cells.Item(102, 2).NumberFormat = "@"    # date column is stored as text
$ws.Cells.Item(102, 2).Value = "2018-11-06"  # date
$ws.Cells.Item(102, 8).Value = "X"           # Senate
$ws.Cells.Item(102, 9).Value = "X"           # House
$ws.Cells.Item(102, 10).Value = "X"          # Governor
$ws.Cells.Item(102, 12).Value = "X"          # Atty General
$ws.Cells.Item(102, 13).Value = "X"          # Sec of St
$ws.Cells.Item(102, 14).Value = "X"          # St Treasurer
$ws.Cells.Item(102, 16).Value = "X"          # St Senate
$ws.Cells.Item(102, 17).Value = "X"          # St Assembly

# Reflect the author's final cursor position / selection from the edit.
$ws.Range("B105").Select()
